# Week8_Recap.pptx — Char slide ("Functions (from ctype.h)") gains a new
# bullet for putchar() right after the existing toupper() bullet, and
# right before the trailing blank sub-bullet paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

# The "toupper()" bullet is the last populated lvl-1 item (paragraph 11);
# paragraph 12 is the trailing empty lvl-1 paragraph that must stay last.
$lastBullet = $tr.Paragraphs(11, 1)

# Insert a new paragraph after it. This clones the paragraph/run
# formatting of the preceding bullet (same pPr: marL/lvl/indent/bullet/tabs).
$newRun = $lastBullet.InsertAfter([char]13 + "putchar()")

# Split "putchar()" into two runs ("putchar" / "()") matching the
# existing isalpha()/isupper()/... bullets, which all separate the
# function name from the trailing parentheses.
$full = $shp.TextFrame.TextRange
$newParagraph = $full.Paragraphs(12, 1)
$parenStart = $newParagraph.Start + $newParagraph.Length - 3
$parens = $full.Characters($parenStart, 2)
$parens.Font.Bold = $true
